$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "%"
$ws.Range("J1").Value = "%Flux"

$ws.Range("I2:I98").Formula = "=E2/D2"
$ws.Range("J2:J98").Formula = "=H2/F2"
